# Add four new vehicles to the Tracking Table (rows 56-59) and fix the
# "D" column typo bug (new rows correctly reference Heavy/Light Goods),
# replicating the pattern already present in the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for the four new vehicles to append, in order.
$newVehicles = @(
    @{ Name = "Sentinel DG4";        Year = 1928; Order = 1; Type = "Heavy Goods"; Top = 37; Cap = 15 },
    @{ Name = "AEC Mammoth Major 8"; Year = 1935; Order = 1; Type = "Heavy Goods"; Top = 35; Cap = 24 },
    @{ Name = "Commer FC";           Year = 1960; Order = 1; Type = "Light Goods"; Top = 60; Cap = 8  },
    @{ Name = "Austin 10hp";         Year = 1938; Order = 1; Type = "Light Goods"; Top = 45; Cap = 5  }
)

$startRow = 56
for ($i = 0; $i -lt $newVehicles.Count; $i++) {
    $row = $startRow + $i
    $v = $newVehicles[$i]

    $ws.Cells.Item($row, 1).Value = $v.Name
    $ws.Cells.Item($row, 2).Value = $v.Year
    $ws.Cells.Item($row, 3).Value = $v.Order
    $ws.Cells.Item($row, 4).Value = $v.Type

    $ws.Cells.Item($row, 5).Formula = "=IF(B$row > 1900, ((B$row-1900)*10)+400+C$row, ((B$row-1730)*2)+C$row)+VLOOKUP(D$row,'ID Scheme'!`$A`$2:`$B`$6,2, FALSE)"

    $ws.Cells.Item($row, 6).Value = $v.Top
    $ws.Cells.Item($row, 7).Value = $v.Cap

    $ws.Cells.Item($row, 8).Formula = "=SQRT(F$row*G$row)/`$B`$1"
    $ws.Cells.Item($row, 8).NumberFormat = $ws.Cells.Item($row - 1, 8).NumberFormat

    $ws.Cells.Item($row, 9).Formula = "=H$row*0.9"
    $ws.Cells.Item($row, 9).NumberFormat = $ws.Cells.Item($row - 1, 9).NumberFormat

    $ws.Cells.Item($row, 10).Value = "x"
    $ws.Cells.Item($row, 10).NumberFormat = $ws.Cells.Item($row - 1, 10).NumberFormat
}

# Mirror the author's final selection / scroll position in the saved view.
$ws.Activate() | Out-Null
$ws.Range("G57").Select() | Out-Null

$wb.Save()
